# Tire Type Filtering for dashboard script and cleanup of Tire Type extraction
# in process_audio_to_csv script.
#
# The underlying per-segment intensity values (Step1_Data) were recomputed
# after the Tire_Type filtering / extraction cleanup, which changes which
# audio segments are attributed to each tire. This script writes the new
# Step1_Data values and then rebuilds the fully-derived downstream sheets
# (Step2_Sj cumulative sums, and the Step3_DataPts_* threshold summaries)
# exactly the way the original pipeline produced them, rather than
# hard-coding every derived literal.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("Step1_Data")
$ws2 = $wb.Worksheets.Item("Step2_Sj")

# --- 1. Updated raw Step1_Data values (new Tire_Type filtering results) ---
# Row 2 (signal segment 1)
$ws1.Range("D2").Value = 0.1811435519263472
$ws1.Range("E2").Value = 0.009316160306502454
$ws1.Range("F2").Value = 0.2010525851117242
$ws1.Range("G2").Value = 0.03145834454140034
$ws1.Range("H2").Value = 0.0859089834946998
$ws1.Range("I2").Value = 0.008095205408709109
$ws1.Range("K2").Value = 0.04550472331419719
$ws1.Range("L2").Value = 0.01330009693577723
$ws1.Range("M2").Value = 0.06774056254240343
$ws1.Range("O2").Value = 0.04360571346702332
$ws1.Range("P2").Value = 0.1937244977383105
$ws1.Range("R2").Value = 0.1051872775479908
$ws1.Range("X2").Value = 0.007797359641352123
$ws1.Range("AD2").Value = 0.006164938023562178

# Row 3 (signal segment 2)
$ws1.Range("D3").Value = 0.08998202977656139
$ws1.Range("E3").Value = 0.1338324908650348
$ws1.Range("F3").Value = 0.3264599260158483
$ws1.Range("G3").Value = 0.02527205564436814
$ws1.Range("H3").Value = 0.02884738546717558
$ws1.Range("I3").Value = 0.006231418202691876
$ws1.Range("M3").Value = 0.04656976059326211
$ws1.Range("N3").Value = 0.003754004880391086
$ws1.Range("O3").Value = 0.006369942228464359
$ws1.Range("P3").Value = 0.1483655530225822
$ws1.Range("Q3").Value = 0.004842568737515367
$ws1.Range("R3").Value = 0.07630433114304561
$ws1.Range("S3").Value = 0.03816950395917588
$ws1.Range("X3").Value = 0.02082600240543147
$ws1.Range("Y3").Value = 0.001665447083112227
$ws1.Range("Z3").Value = 0.0004923665354176646
$ws1.Range("AA3").Value = 0.02427883158089236
$ws1.Range("AD3").Value = 0.01773638185902971

# Row 4 (signal segment 3)
$ws1.Range("D4").Value = 0.003583621028066395
$ws1.Range("E4").Value = 0.2587782657676923
$ws1.Range("F4").Value = 0.1425581228677392
$ws1.Range("G4").Value = 0.07109037717135615
$ws1.Range("H4").Value = 0.04485901163006848
$ws1.Range("I4").Value = 0.02864700980879237
$ws1.Range("K4").Value = 0.05497789896731796
$ws1.Range("L4").Value = 0.005755054952078799
$ws1.Range("M4").Value = 0.03053889852826559
$ws1.Range("N4").Value = 0.01959962402452297
$ws1.Range("P4").Value = 0.1835952202584066
$ws1.Range("R4").Value = 0.1540364486225093
$ws1.Range("T4").Value = 0.001980446373183933

# Row 5 (signal segment 4)
$ws1.Range("D5").Value = 0.2100262200267631
$ws1.Range("F5").Value = 0.1831355023249108
$ws1.Range("G5").Value = 0.05170601472930404
$ws1.Range("H5").Value = 0.07531233250480017
$ws1.Range("I5").Value = 0.001215943094987203
$ws1.Range("K5").Value = 0.02943893810326729
$ws1.Range("L5").Value = 0.009310019696327563
$ws1.Range("M5").Value = 0.08013640094626159
$ws1.Range("O5").Value = 0.04278807937289045
$ws1.Range("P5").Value = 0.188164921075956
$ws1.Range("R5").Value = 0.1212251315524002
$ws1.Range("X5").Value = 0.003315982381709649
$ws1.Range("AD5").Value = 0.004224514190421797

# Row 6 (signal segment 5)
$ws1.Range("E6").Value = 0.2659551940545829
$ws1.Range("F6").Value = 0.1583131664833482
$ws1.Range("G6").Value = 0.07690656275541022
$ws1.Range("H6").Value = 0.05525486747094241
$ws1.Range("I6").Value = 0.0253837366350721
$ws1.Range("J6").Value = 0.006591238131961781
$ws1.Range("K6").Value = 0.0458310038727204
$ws1.Range("M6").Value = 0.01051483876344519
$ws1.Range("N6").Value = 0.03683909674025296
$ws1.Range("P6").Value = 0.1571572767149673
$ws1.Range("Q6").Value = 0.01327963174575682
$ws1.Range("R6").Value = 0.1258400449456436
$ws1.Range("T6").Value = 0.01035362388352389
$ws1.Range("AD6").Value = 0.0117797178023723

# --- 2. Rebuild Step2_Sj: running cumulative sum across columns B..AJ per row ---
$firstCol = 2   # column B
$lastCol = 36   # column AJ

for ($r = 2; $r -le 6; $r++) {
    $cum = 0.0
    for ($c = $firstCol; $c -le $lastCol; $c++) {
        $cum = $cum + $ws1.Cells.Item($r, $c).Value()
        $ws2.Cells.Item($r, $c).Value = $cum
    }
}

# --- 3. Rebuild the Step3_DataPts_* threshold summary sheets from Step2_Sj ---
# Columns: B=Intensity_Threshold, C=Tire_Number(unchanged), D=First_Noticeable_Increase_Index,
#          E=Point_Exceeds_Index(unused/0), F=First_Noticeable_Increase_Cumulative_Value,
#          G=Point_Exceeds_Cumulative_Value(=D-C)
$thresholdSheets = @(
    @{ Name = "Step3_DataPts_0.5"; Threshold = 0.5 },
    @{ Name = "Step3_DataPts_0.7"; Threshold = 0.7 },
    @{ Name = "Step3_DataPts_0.8"; Threshold = 0.8 },
    @{ Name = "Step3_DataPts_0.9"; Threshold = 0.9 }
)

foreach ($sheetInfo in $thresholdSheets) {
    $ws3 = $wb.Worksheets.Item($sheetInfo.Name)
    $threshold = $sheetInfo.Threshold

    for ($r = 2; $r -le 6; $r++) {
        $tireNumber = $ws3.Cells.Item($r, 3).Value()   # column C, unchanged

        $foundCol = -1
        $foundVal = 0.0
        for ($c = $firstCol; $c -le $lastCol; $c++) {
            $val = $ws2.Cells.Item($r, $c).Value()
            if ($foundCol -eq -1 -and $val -gt $threshold) {
                $foundCol = $c
                $foundVal = $val
            }
        }

        $index = $foundCol - $firstCol + 1   # 1-based index among B..AJ

        $ws3.Cells.Item($r, 4).Value = $index                  # D
        $ws3.Cells.Item($r, 6).Value = $foundVal                # F
        $ws3.Cells.Item($r, 7).Value = ($index - $tireNumber)   # G
    }
}
